# Auto-generated PowerShell COM-interop script
# Applies numeric updates to the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
# per the scheduled-runner price refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(74, 8).Value = 3030.4  # H74: was 3478.2222
$ws.Cells.Item(74, 9).Value = 3000.5  # I74: was 3333.8333
$ws.Cells.Item(74, 10).Value = 3075.25  # J74: was 3767
$ws.Cells.Item(74, 11).Value = 3000.5  # K74: was 3333.8333
$ws.Cells.Item(74, 12).Value = 3075.25  # L74: was 3767
$ws.Cells.Item(74, 13).Value = -2064.5  # M74: was -2397.8333
$ws.Cells.Item(74, 14).Value = -4947.25  # N74: was -5639

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(77, 8).Value = 3030.4  # H77: was 3478.2222
$ws.Cells.Item(77, 9).Value = 3000.5  # I77: was 3333.8333
$ws.Cells.Item(77, 10).Value = 3075.25  # J77: was 3767
$ws.Cells.Item(77, 11).Value = 15002.5  # K77: was 16669.1665
$ws.Cells.Item(77, 12).Value = 15376.25  # L77: was 18835
$ws.Cells.Item(77, 13).Value = -10322.5  # M77: was -11989.1665
$ws.Cells.Item(77, 14).Value = -24736.25  # N77: was -28195

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(88, 8).Value = 3086795.5  # H88: was 2058461.1
$ws.Cells.Item(88, 9).Value = 301.5  # I88: was 396.66666
$ws.Cells.Item(88, 10).Value = 6173289.5  # J88: was 4116525.8
$ws.Cells.Item(88, 11).Value = 301.5  # K88: was 396.66666
$ws.Cells.Item(88, 12).Value = 6173289.5  # L88: was 4116525.8
$ws.Cells.Item(88, 13).Value = 104.5  # M88: was 9.333340000000021
$ws.Cells.Item(88, 14).Value = -6174101.5  # N88: was -4117337.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(91, 8).Value = 3086795.5  # H91: was 2058461.1
$ws.Cells.Item(91, 9).Value = 301.5  # I91: was 396.66666
$ws.Cells.Item(91, 10).Value = 6173289.5  # J91: was 4116525.8
$ws.Cells.Item(91, 11).Value = 301.5  # K91: was 396.66666
$ws.Cells.Item(91, 12).Value = 6173289.5  # L91: was 4116525.8
$ws.Cells.Item(91, 13).Value = 1102.5  # M91: was 1007.33334
$ws.Cells.Item(91, 14).Value = -6176097.5  # N91: was -4119333.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(97, 8).Value = 2087.318  # H97: was 2286.8
$ws.Cells.Item(97, 10).Value = 2087.318  # J97: was 2286.8
$ws.Cells.Item(97, 12).Value = 6261.954000000001  # L97: was 6860.400000000001
$ws.Cells.Item(97, 14).Value = -7253.954000000001  # N97: was -7852.400000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(100, 8).Value = 2463  # H100: was 2719
$ws.Cells.Item(100, 10).Value = 0  # J100: was 3999
$ws.Cells.Item(100, 12).Value = 0  # L100: was 3999
$ws.Cells.Item(100, 14).ClearContents()  # N100: was -5081

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(112, 8).Value = 2075.15  # H112: was 2700.2727
$ws.Cells.Item(112, 9).Value = 671  # I112: was 649.625
$ws.Cells.Item(112, 10).Value = 2373  # J112: was 3872.0715
$ws.Cells.Item(112, 11).Value = 2013  # K112: was 1948.875
$ws.Cells.Item(112, 12).Value = 7119  # L112: was 11616.2145
$ws.Cells.Item(112, 13).Value = -905  # M112: was -840.875
$ws.Cells.Item(112, 14).Value = -9335  # N112: was -13832.2145

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(131, 8).Value = 0  # H131: was 1763.3334
$ws.Cells.Item(131, 9).Value = 0  # I131: was 2145
$ws.Cells.Item(131, 10).Value = 0  # J131: was 1000
$ws.Cells.Item(131, 11).Value = 0  # K131: was 6435
$ws.Cells.Item(131, 12).Value = 0  # L131: was 3000
$ws.Cells.Item(131, 13).ClearContents()  # M131: was -1395
$ws.Cells.Item(131, 14).ClearContents()  # N131: was -13080

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 9808107  # H132: was 4674.8335
$ws.Cells.Item(132, 9).Value = 13891087  # I132: was 2343.4092
$ws.Cells.Item(132, 10).Value = 8954  # J132: was 11086.25
$ws.Cells.Item(132, 11).Value = 41673261  # K132: was 7030.2276
$ws.Cells.Item(132, 12).Value = 26862  # L132: was 33258.75
$ws.Cells.Item(132, 13).Value = -41670731  # M132: was -4500.2276
$ws.Cells.Item(132, 14).Value = -31922  # N132: was -38318.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 2189.94  # H138: was 822152.5
$ws.Cells.Item(138, 9).Value = 1423  # I138: was 1102.8948
$ws.Cells.Item(138, 10).Value = 2294.5227  # J138: was 1280974.4
$ws.Cells.Item(138, 11).Value = 4269  # K138: was 3308.6844
$ws.Cells.Item(138, 12).Value = 6883.5681  # L138: was 3842923.2
$ws.Cells.Item(138, 13).Value = 871  # M138: was 1831.3156
$ws.Cells.Item(138, 14).Value = -17163.5681  # N138: was -3853203.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2051.5  # H2: was 1884.1666
$ws.Cells.Item(2, 9).Value = 954.7778  # I2: was 958.4
$ws.Cells.Item(2, 10).Value = 5341.6665  # J2: was 6513
$ws.Cells.Item(2, 11).Value = 954.7778  # K2: was 958.4
$ws.Cells.Item(2, 12).Value = 5341.6665  # L2: was 6513
$ws.Cells.Item(2, 13).Value = -841.7778  # M2: was -845.4
$ws.Cells.Item(2, 14).Value = -5567.6665  # N2: was -6739

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3079.413  # H32: was 3327.762
$ws.Cells.Item(32, 9).Value = 3246.372  # I32: was 3459.025
$ws.Cells.Item(32, 10).Value = 686.3333  # J32: was 702.5
$ws.Cells.Item(32, 11).Value = 3246.372  # K32: was 3459.025
$ws.Cells.Item(32, 12).Value = 686.3333  # L32: was 702.5
$ws.Cells.Item(32, 13).Value = -2959.372  # M32: was -3172.025
$ws.Cells.Item(32, 14).Value = -1260.3333  # N32: was -1276.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(53, 8).Value = 4750  # H53: was 6000
$ws.Cells.Item(53, 9).Value = 2333.3333  # I53: was 3000
$ws.Cells.Item(53, 11).Value = 2333.3333  # K53: was 3000
$ws.Cells.Item(53, 13).Value = -1651.3333  # M53: was -2318

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(116, 8).Value = 2051.5  # H116: was 1884.1666
$ws.Cells.Item(116, 9).Value = 954.7778  # I116: was 958.4
$ws.Cells.Item(116, 10).Value = 5341.6665  # J116: was 6513
$ws.Cells.Item(116, 11).Value = 954.7778  # K116: was 958.4
$ws.Cells.Item(116, 12).Value = 5341.6665  # L116: was 6513
$ws.Cells.Item(116, 13).Value = 1339.2222  # M116: was 1335.6
$ws.Cells.Item(116, 14).Value = -9929.666499999999  # N116: was -11101

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 1363  # H122: was 1489.75
$ws.Cells.Item(122, 9).Value = 1116.5714  # I122: was 1220.8
$ws.Cells.Item(122, 11).Value = 3349.7142  # K122: was 3662.4
$ws.Cells.Item(122, 13).Value = -899.7142000000003  # M122: was -1212.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2051.5  # H3: was 1884.1666
$ws.Cells.Item(3, 9).Value = 954.7778  # I3: was 958.4
$ws.Cells.Item(3, 10).Value = 5341.6665  # J3: was 6513
$ws.Cells.Item(3, 11).Value = 954.7778  # K3: was 958.4
$ws.Cells.Item(3, 12).Value = 5341.6665  # L3: was 6513
$ws.Cells.Item(3, 13).Value = -840.7778  # M3: was -844.4
$ws.Cells.Item(3, 14).Value = -5569.6665  # N3: was -6741

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(87, 8).Value = 75000  # H87: was 59666.668
$ws.Cells.Item(87, 10).Value = 75000  # J87: was 59666.668
$ws.Cells.Item(87, 12).Value = 75000  # L87: was 59666.668
$ws.Cells.Item(87, 14).Value = -77496  # N87: was -62162.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(90, 8).Value = 75000  # H90: was 59666.668
$ws.Cells.Item(90, 10).Value = 75000  # J90: was 59666.668
$ws.Cells.Item(90, 12).Value = 225000  # L90: was 179000.004
$ws.Cells.Item(90, 14).Value = -237480  # N90: was -191480.004

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 946.5294  # H58: was 744.4167
$ws.Cells.Item(58, 9).Value = 939.4  # I58: was 726.8570999999999
$ws.Cells.Item(58, 10).Value = 1000  # J58: was 769
$ws.Cells.Item(58, 11).Value = 939.4  # K58: was 726.8570999999999
$ws.Cells.Item(58, 12).Value = 1000  # L58: was 769
$ws.Cells.Item(58, 13).Value = -736.4  # M58: was -523.8570999999999
$ws.Cells.Item(58, 14).Value = -1406  # N58: was -1175

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 2020.5714  # H99: was 1998.2858
$ws.Cells.Item(99, 9).Value = 1826  # I99: was 1794.8
$ws.Cells.Item(99, 11).Value = 1826  # K99: was 1794.8
$ws.Cells.Item(99, 13).Value = -328  # M99: was -296.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(126, 8).Value = 2020.5714  # H126: was 1998.2858
$ws.Cells.Item(126, 9).Value = 1826  # I126: was 1794.8
$ws.Cells.Item(126, 11).Value = 5478  # K126: was 5384.4
$ws.Cells.Item(126, 13).Value = -3008  # M126: was -2914.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 946.5294  # H136: was 744.4167
$ws.Cells.Item(136, 9).Value = 939.4  # I136: was 726.8570999999999
$ws.Cells.Item(136, 10).Value = 1000  # J136: was 769
$ws.Cells.Item(136, 11).Value = 2818.2  # K136: was 2180.5713
$ws.Cells.Item(136, 12).Value = 3000  # L136: was 2307
$ws.Cells.Item(136, 13).Value = -268.1999999999998  # M136: was 369.4287000000004
$ws.Cells.Item(136, 14).Value = -8100  # N136: was -7407

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(140, 8).Value = 30900  # H140: was 40000
$ws.Cells.Item(140, 10).Value = 30900  # J140: was 40000
$ws.Cells.Item(140, 12).Value = 30900  # L140: was 40000
$ws.Cells.Item(140, 14).Value = -41260  # N140: was -50360

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(141, 8).Value = 33936  # H141: was 34000
$ws.Cells.Item(141, 10).Value = 33926.855  # J141: was 34000
$ws.Cells.Item(141, 12).Value = 33926.855  # L141: was 34000
$ws.Cells.Item(141, 14).Value = -44286.855  # N141: was -44360

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 655.8  # H12: was 796.6667
$ws.Cells.Item(12, 9).Value = 439  # I12: was 0
$ws.Cells.Item(12, 10).Value = 679.8889  # J12: was 796.6667
$ws.Cells.Item(12, 11).Value = 1317  # K12: was 0
$ws.Cells.Item(12, 12).Value = 2039.6667  # L12: was 2390.0001
$ws.Cells.Item(12, 13).Value = -1144  # M12: new cell
$ws.Cells.Item(12, 14).Value = -2385.6667  # N12: was -2736.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 5771.3  # H107: was 5082
$ws.Cells.Item(107, 10).Value = 13497.375  # J107: was 9949
$ws.Cells.Item(107, 12).Value = 40492.125  # L107: was 29847
$ws.Cells.Item(107, 14).Value = -44332.125  # N107: was -33687

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 922  # H122: was 932.1667
$ws.Cells.Item(122, 9).Value = 0  # I122: was 572.25
$ws.Cells.Item(122, 10).Value = 922  # J122: was 1112.125
$ws.Cells.Item(122, 11).Value = 0  # K122: was 5150.25
$ws.Cells.Item(122, 12).Value = 8298  # L122: was 10009.125
$ws.Cells.Item(122, 13).ClearContents()  # M122: was -2700.25
$ws.Cells.Item(122, 14).Value = -13198  # N122: was -14909.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(134, 8).Value = 3191.8096  # H134: was 3581.5557
$ws.Cells.Item(134, 9).Value = 1862.9  # I134: was 2295.5715
$ws.Cells.Item(134, 11).Value = 5588.700000000001  # K134: was 6886.7145
$ws.Cells.Item(134, 13).Value = -518.7000000000007  # M134: was -1816.7145

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 593.4375  # H107: was 606.0714
$ws.Cells.Item(107, 9).Value = 412.27274  # I107: was 452.5
$ws.Cells.Item(107, 10).Value = 992  # J107: was 990
$ws.Cells.Item(107, 11).Value = 412.27274  # K107: was 452.5
$ws.Cells.Item(107, 12).Value = 992  # L107: was 990
$ws.Cells.Item(107, 13).Value = 1507.72726  # M107: was 1467.5
$ws.Cells.Item(107, 14).Value = -4832  # N107: was -4830

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2373.1936  # H132: was 2430.0688
$ws.Cells.Item(132, 9).Value = 1953.7  # I132: was 2009.2632
$ws.Cells.Item(132, 10).Value = 3135.9092  # J132: was 3229.6
$ws.Cells.Item(132, 11).Value = 5861.1  # K132: was 6027.7896
$ws.Cells.Item(132, 12).Value = 9407.7276  # L132: was 9688.799999999999
$ws.Cells.Item(132, 13).Value = -3331.1  # M132: was -3497.7896
$ws.Cells.Item(132, 14).Value = -14467.7276  # N132: was -14748.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 1225.25  # H46: was 882
$ws.Cells.Item(46, 9).Value = 600.5  # I46: was 330.5
$ws.Cells.Item(46, 11).Value = 600.5  # K46: was 330.5
$ws.Cells.Item(46, 13).Value = -412.5  # M46: was -142.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 2559.7273  # H136: was 2137.0667
$ws.Cells.Item(136, 9).Value = 2819.125  # I136: was 2616.889
$ws.Cells.Item(136, 10).Value = 1868  # J136: was 1417.3334
$ws.Cells.Item(136, 11).Value = 8457.375  # K136: was 7850.667
$ws.Cells.Item(136, 12).Value = 5604  # L136: was 4252.0002
$ws.Cells.Item(136, 13).Value = -5907.375  # M136: was -5300.667
$ws.Cells.Item(136, 14).Value = -10704  # N136: was -9352.0002

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 333335680  # H2: was 5000
$ws.Cells.Item(2, 9).Value = 500002500  # I2: was 5000
$ws.Cells.Item(2, 10).Value = 2000  # J2: was 0
$ws.Cells.Item(2, 11).Value = 500002500  # K2: was 5000
$ws.Cells.Item(2, 12).Value = 2000  # L2: was 0
$ws.Cells.Item(2, 13).Value = -500002388  # M2: was -4888
$ws.Cells.Item(2, 14).Value = -2224  # N2: new cell

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(41, 8).Value = 10342  # H41: was 9100.75
$ws.Cells.Item(41, 9).Value = 10342  # I41: was 8675.333000000001
$ws.Cells.Item(41, 10).Value = 0  # J41: was 10377
$ws.Cells.Item(41, 11).Value = 10342  # K41: was 8675.333000000001
$ws.Cells.Item(41, 12).Value = 0  # L41: was 10377
$ws.Cells.Item(41, 13).ClearContents()  # M41: was -8285.333000000001
$ws.Cells.Item(41, 14).Value = -9952  # N41: was -11157

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 16252119  # H122: was 17335520
$ws.Cells.Item(122, 9).Value = 17335394  # I122: was 18573556
$ws.Cells.Item(122, 11).Value = 52006182  # K122: was 55720668
$ws.Cells.Item(122, 13).Value = -52003732  # M122: was -55718218
